$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.096.22"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.476.40"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'585.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'131.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("E9").Value = "  +5.87%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "4.068.51"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "3.474.15"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "64.084.79"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'385.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "3.615.54"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'74.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "'7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").Value = "'7.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("E32").Value = "  -4.42%  "
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("D34").Value = "3.503.02"
$ws.Range("D36").Value = "'22.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").Value = "'5.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("D40").Value = "'161.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'41.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "'1.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "'23.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.92%  "
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "2.331.01"
$ws.Range("E51").Value = "  -5.25%  "
